$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Urban"/"Город"/"Шаар" row (row 23) to the new, more specific wording.
$ws.Range("A23").Value = "Шаар жерлери"
$ws.Range("B23").Value = "Городские поселения"
$ws.Range("C23").Value = "City"

# Update the "Rural"/"Село"/"Айыл" row (row 24) to the new, more specific wording.
$ws.Range("A24").Value = "Айыл аймагы"
$ws.Range("B24").Value = "Сельская местность"
$ws.Range("C24").Value = "Village"

# Update the saved cursor/selection position.
$ws.Range("C30").Select()
